$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.292394999999999
$ws.Range("H2").Value = 21.877185
$ws.Range("I2").Value = 0.5244715940033005
$ws.Range("J2").Value = 0.5244715940033005
$ws.Range("M2").Value = 7.292394999999999
$ws.Range("N2").Value = 21.877185
$ws.Range("O2").Value = 0.5244715940033005
$ws.Range("P2").Value = 0.5244715940033005
$ws.Range("Q2").Value = 53.17902483602499
$ws.Range("R2").Value = 478.6112235242249
$ws.Range("S2").Value = 0.2750704529163628
$ws.Range("T2").Value = 0.2750704529163628
$ws.Range("G3").Value = 7.292394999999999
$ws.Range("H3").Value = 21.877185
$ws.Range("I3").Value = 0.5244715940033005
$ws.Range("J3").Value = 0.5244715940033005
$ws.Range("O3").Value = 0.02354183170388992
$ws.Range("P3").Value = 0.02354183170388992
$ws.Range("Q3").Value = 2.38703424014
$ws.Range("R3").Value = 21.48330816126
$ws.Range("S3").Value = 0.01234702199949658
$ws.Range("T3").Value = 0.01234702199949658
$ws.Range("G4").Value = 7.292394999999999
$ws.Range("H4").Value = 21.877185
$ws.Range("I4").Value = 0.5244715940033005
$ws.Range("J4").Value = 0.5244715940033005
$ws.Range("O4").Value = 0.4519865742928097
$ws.Range("P4").Value = 0.4519865742928096
$ws.Range("Q4").Value = 45.82937481208166
$ws.Range("R4").Value = 412.4643733087349
$ws.Range("S4").Value = 0.2370541190874411
$ws.Range("T4").Value = 0.237054119087441
$ws.Range("I5").Value = 0.02354183170388992
$ws.Range("J5").Value = 0.02354183170388992
$ws.Range("M5").Value = 7.292394999999999
$ws.Range("N5").Value = 21.877185
$ws.Range("O5").Value = 0.5244715940033005
$ws.Range("P5").Value = 0.5244715940033005
$ws.Range("Q5").Value = 2.38703424014
$ws.Range("R5").Value = 21.48330816126
$ws.Range("S5").Value = 0.01234702199949658
$ws.Range("T5").Value = 0.01234702199949658
$ws.Range("I6").Value = 0.02354183170388992
$ws.Range("J6").Value = 0.02354183170388992
$ws.Range("O6").Value = 0.02354183170388992
$ws.Range("P6").Value = 0.02354183170388992
$ws.Range("S6").Value = 0.0005542178399742768
$ws.Range("T6").Value = 0.0005542178399742768
$ws.Range("I7").Value = 0.02354183170388992
$ws.Range("J7").Value = 0.02354183170388992
$ws.Range("O7").Value = 0.4519865742928097
$ws.Range("P7").Value = 0.4519865742928096
$ws.Range("S7").Value = 0.01064059186441907
$ws.Range("T7").Value = 0.01064059186441906
$ws.Range("I8").Value = 0.4519865742928097
$ws.Range("J8").Value = 0.4519865742928096
$ws.Range("M8").Value = 7.292394999999999
$ws.Range("N8").Value = 21.877185
$ws.Range("O8").Value = 0.5244715940033005
$ws.Range("P8").Value = 0.5244715940033005
$ws.Range("Q8").Value = 45.82937481208166
$ws.Range("R8").Value = 412.4643733087349
$ws.Range("S8").Value = 0.2370541190874411
$ws.Range("T8").Value = 0.237054119087441
$ws.Range("I9").Value = 0.4519865742928097
$ws.Range("J9").Value = 0.4519865742928096
$ws.Range("O9").Value = 0.02354183170388992
$ws.Range("P9").Value = 0.02354183170388992
$ws.Range("S9").Value = 0.01064059186441907
$ws.Range("T9").Value = 0.01064059186441906
$ws.Range("I10").Value = 0.4519865742928097
$ws.Range("J10").Value = 0.4519865742928096
$ws.Range("O10").Value = 0.4519865742928097
$ws.Range("P10").Value = 0.4519865742928096
$ws.Range("S10").Value = 0.2042918633409496
$ws.Range("T10").Value = 0.2042918633409495
